$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff / Handback Datetime for row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-20 05:09:35"
$wsZhCn.Range("H2").Value = "2016-03-20 05:10:17"

# de-de sheet: update Correspond Handoff / Handback Datetime for row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-20 05:09:43"
$wsDeDe.Range("H2").Value = "2016-03-20 05:10:33"
